$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T3").Value = 110
$ws.Range("T4").Value = 131
$ws.Range("T5").Value = 70
$ws.Range("T6").Value = 108
$ws.Range("T7").Value = 188
$ws.Range("T8").Value = 87
$ws.Range("T9").Value = 151
$ws.Range("T10").Value = 206
$ws.Range("T11").Value = 132
$ws.Range("T12").Value = 105
$ws.Range("T13").Value = 132
$ws.Range("T14").Value = 178
$ws.Range("T15").Value = 44
$ws.Range("T16").Value = 204
$ws.Range("T17").Value = 197
$ws.Range("T18").Value = 38
$ws.Range("T19").Value = 82
$ws.Range("T20").Value = 109
$ws.Range("T21").Value = 96
$ws.Range("T22").Value = 91
$ws.Range("T23").Value = 215
$ws.Range("T24").Value = 206
$ws.Range("T25").Value = 62
$ws.Range("T26").Value = 75
$ws.Range("T27").Value = 125
$ws.Range("T28").Value = 73
$ws.Range("T29").Value = 72
$ws.Range("T30").Value = 61
$ws.Range("T31").Value = 39
$ws.Range("T32").Value = 34
$ws.Range("T33").Value = 25
$ws.Range("T34").Value = 38
